$wb = $excel.ActiveWorkbook

# --- Technologies sheet: remove ".Net" row content, fix last row to "CSS" ---
$wsTech = $wb.Worksheets.Item("Technologies")
$wsTech.Range("B3").Value = ""
$wsTech.Range("B7").Value = "CSS"

# --- Employees sheet: correct duplicated EMP ID E0128 -> E0127 ---
$wsEmp = $wb.Worksheets.Item("Employees")
$wsEmp.Range("B9").Value = "E0127"

# --- Roles sheet: remove completed / duplicate role rows ---
$wsRoles = $wb.Worksheets.Item("Roles")
$wsRoles.Range("B3").Value = ""
$wsRoles.Range("B5").Value = "Data engineer 4"
$wsRoles.Range("B6").Value = ""

# --- Selection / active sheet bookkeeping (matches author's final click state) ---
$wsRoles.Range("H14").Select()
$wsEmp.Range("B9").Select()
$wsTech.Activate()
$wsTech.Range("B7").Select()
